$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '44.473.26'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +3.62%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.425.84'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +2.79%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.07%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '313.20'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +3.69%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '101.66'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +6.50%  '
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +2.03%  '
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -0.08%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.514'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +5.87%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '35.21'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +4.24%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0799'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +2.10%  '
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +1.48%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '18.87'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +3.01%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.94'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +3.51%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '2.804.52'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +2.78%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.453.12'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +4.54%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.836'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +5.37%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '44.383.57'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +3.50%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.40'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +4.98%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.39'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +2.44%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.0#0908'
$ws.Range("D21").Replace('#', [char]0x2083) | Out-Null
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +2.89%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '68.89'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +1.52%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '240.65'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +2.49%  '
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +4.64%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.47'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +2.27%  '
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +0.03%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '25.15'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +2.48%  '
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -4.20%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.61'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +4.21%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '33.26'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +5.63%  '
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +1.14%  '
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +17.52%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '19.54'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +13.32%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.17'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +3.39%  '
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +0.23%  '
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +7.05%  '
$ws.Range("B37").NumberFormat = "@"
$ws.Range("B37").Value = 'ARBITRUM'
$ws.Range("C37").NumberFormat = "@"
$ws.Range("C37").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.90'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +3.41%  '
$ws.Range("B38").NumberFormat = "@"
$ws.Range("B38").Value = 'RenderToken'
$ws.Range("C38").NumberFormat = "@"
$ws.Range("C38").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '4.52'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +4.48%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.89'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +4.37%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '125.44'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +1.74%  '
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +1.32%  '
$ws.Range("B42").NumberFormat = "@"
$ws.Range("B42").Value = 'WEMIXToken'
$ws.Range("C42").NumberFormat = "@"
$ws.Range("C42").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.16'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -5.80%  '
$ws.Range("B43").NumberFormat = "@"
$ws.Range("B43").Value = 'EnergySwap'
$ws.Range("C43").NumberFormat = "@"
$ws.Range("C43").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '21.46'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +2.38%  '
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +3.84%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.945.41'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +0.70%  '
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +1.72%  '
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +9.12%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '9.60'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +5.34%  '
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +10.76%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '53.49'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +4.03%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '73.53'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +3.07%  '
